$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing row (19) down into the new row 20,
# then overwrite with the new row's values.
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.043309689777173
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 1.002299702378884
